$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2131.9707
$ws.Range("I15").Value = 2131.9707
$ws.Range("K15").Value = 6395.9121
$ws.Range("M15").Value = -6226.9121

$ws.Range("H33").Value = 2949.3333
$ws.Range("I33").Value = 3336.3547
$ws.Range("J33").Value = 549.8
$ws.Range("K33").Value = 3336.3547
$ws.Range("L33").Value = 549.8
$ws.Range("M33").Value = -3107.3547
$ws.Range("N33").Value = -1007.8

$ws.Range("H76").Value = 7268.75
$ws.Range("I76").Value = 5525.1665
$ws.Range("K76").Value = 5525.1665
$ws.Range("M76").Value = -5210.1665

$ws.Range("H79").Value = 7268.75
$ws.Range("I79").Value = 5525.1665
$ws.Range("K79").Value = 5525.1665
$ws.Range("M79").Value = -4433.1665

$ws.Range("H80").Value = 1509.6471
$ws.Range("I80").Value = 1346.375
$ws.Range("J80").Value = 1654.7778
$ws.Range("K80").Value = 4039.125
$ws.Range("L80").Value = 4964.3334
$ws.Range("M80").Value = -3041.125
$ws.Range("N80").Value = -6960.3334

$ws.Range("H83").Value = 1509.6471
$ws.Range("I83").Value = 1346.375
$ws.Range("J83").Value = 1654.7778
$ws.Range("K83").Value = 12117.375
$ws.Range("L83").Value = 14893.0002
$ws.Range("M83").Value = -7125.375
$ws.Range("N83").Value = -24877.0002

$ws.Range("H87").Value = 90900
$ws.Range("J87").Value = 90900
$ws.Range("L87").Value = 90900
$ws.Range("N87").Value = -93396

$ws.Range("H90").Value = 90900
$ws.Range("J90").Value = 90900
$ws.Range("L90").Value = 272700
$ws.Range("N90").Value = -285180

$ws.Range("H106").Value = 3240.9412
$ws.Range("I106").Value = 3240.9412
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3240.9412
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2609.9412
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 3592.5925
$ws.Range("I132").Value = 2729.7917
$ws.Range("K132").Value = 8189.375100000001
$ws.Range("M132").Value = -5659.375100000001

$ws.Range("H135").Value = 499.6
$ws.Range("I135").Value = 408.72726
$ws.Range("K135").Value = 3678.54534
$ws.Range("M135").Value = -1143.54534

$ws.Range("H138").Value = 2433.2036
$ws.Range("I138").Value = 2127.1072
$ws.Range("J138").Value = 2762.8462
$ws.Range("K138").Value = 6381.321599999999
$ws.Range("L138").Value = 8288.5386
$ws.Range("M138").Value = -1241.321599999999
$ws.Range("N138").Value = -18568.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 39950
$ws.Range("J54").Value = 39950
$ws.Range("L54").Value = 39950
$ws.Range("N54").Value = -41488

$ws.Range("H61").Value = 6646.25
$ws.Range("I61").Value = 6646.25
$ws.Range("K61").Value = 6646.25
$ws.Range("M61").Value = -6434.25

$ws.Range("H74").Value = 1763.3636
$ws.Range("I74").Value = 1713.862
$ws.Range("J74").Value = 2122.25
$ws.Range("K74").Value = 1713.862
$ws.Range("L74").Value = 2122.25
$ws.Range("M74").Value = -839.8620000000001
$ws.Range("N74").Value = -3870.25

$ws.Range("H77").Value = 1763.3636
$ws.Range("I77").Value = 1713.862
$ws.Range("J77").Value = 2122.25
$ws.Range("K77").Value = 8569.310000000001
$ws.Range("L77").Value = 10611.25
$ws.Range("M77").Value = -4201.310000000001
$ws.Range("N77").Value = -19347.25

$ws.Range("H122").Value = 1738.4286
$ws.Range("I122").Value = 1700.4
$ws.Range("K122").Value = 5101.200000000001
$ws.Range("M122").Value = -2651.200000000001

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 20101.56
$ws.Range("I132").Value = 24748.8
$ws.Range("K132").Value = 74246.39999999999
$ws.Range("M132").Value = -71716.39999999999

$ws.Range("H136").Value = 6646.25
$ws.Range("I136").Value = 6646.25
$ws.Range("K136").Value = 19938.75
$ws.Range("M136").Value = -17388.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3506.6667
$ws.Range("I5").Value = 5010
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 5010
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -4897
$ws.Range("N5").Value = -726

$ws.Range("H20").Value = 5557527.5
$ws.Range("I20").Value = 7693638.5
$ws.Range("J20").Value = 3637.6
$ws.Range("K20").Value = 7693638.5
$ws.Range("L20").Value = 3637.6
$ws.Range("M20").Value = -7693391.5
$ws.Range("N20").Value = -4131.6

$ws.Range("H134").Value = 2484.4866
$ws.Range("I134").Value = 2457.8857
$ws.Range("K134").Value = 7373.657099999999
$ws.Range("M134").Value = -4838.657099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1349.25
$ws.Range("I6").Value = 1499
$ws.Range("J6").Value = 1199.5
$ws.Range("K6").Value = 1499
$ws.Range("L6").Value = 1199.5
$ws.Range("M6").Value = -1386
$ws.Range("N6").Value = -1425.5

$ws.Range("H58").Value = 26412.293
$ws.Range("I58").Value = 33132.625
$ws.Range("J58").Value = 2517.7778
$ws.Range("K58").Value = 33132.625
$ws.Range("L58").Value = 2517.7778
$ws.Range("M58").Value = -32929.625
$ws.Range("N58").Value = -2923.7778

$ws.Range("H125").Value = 36508.168
$ws.Range("J125").Value = 36508.168
$ws.Range("L125").Value = 36508.168
$ws.Range("N125").Value = -41428.168

$ws.Range("H134").Value = 60833.47
$ws.Range("I134").Value = 60833.47
$ws.Range("K134").Value = 182500.41
$ws.Range("M134").Value = -179965.41

$ws.Range("H136").Value = 26412.293
$ws.Range("I136").Value = 33132.625
$ws.Range("J136").Value = 2517.7778
$ws.Range("K136").Value = 99397.875
$ws.Range("L136").Value = 7553.3334
$ws.Range("M136").Value = -96847.875
$ws.Range("N136").Value = -12653.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 524.8
$ws.Range("I36").Value = 524.8
$ws.Range("K36").Value = 1574.4
$ws.Range("M36").Value = -1405.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3995.25
$ws.Range("I43").Value = 3995.25
$ws.Range("K43").Value = 3995.25
$ws.Range("M43").Value = -3844.25

$ws.Range("H58").Value = 22776
$ws.Range("J58").Value = 27035
$ws.Range("L58").Value = 27035
$ws.Range("N58").Value = -27589

$ws.Range("H95").Value = 36500
$ws.Range("J95").Value = 36500
$ws.Range("L95").Value = 36500
$ws.Range("N95").Value = -41992

$ws.Range("H102").Value = 3297.4092
$ws.Range("I102").Value = 3339.1428
$ws.Range("J102").Value = 3224.375
$ws.Range("K102").Value = 3339.1428
$ws.Range("L102").Value = 3224.375
$ws.Range("M102").Value = -1717.1428
$ws.Range("N102").Value = -6468.375

$ws.Range("H122").Value = 3845.25
$ws.Range("J122").Value = 5497.5
$ws.Range("L122").Value = 16492.5
$ws.Range("N122").Value = -21392.5

$ws.Range("H132").Value = 28453.615
$ws.Range("I132").Value = 41230.883
$ws.Range("K132").Value = 123692.649
$ws.Range("M132").Value = -121162.649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3999.25
$ws.Range("I7").Value = 3932.3333
$ws.Range("K7").Value = 3932.3333
$ws.Range("M7").Value = -3820.3333

$ws.Range("H22").Value = 59749.316
$ws.Range("I22").Value = 123793.89
$ws.Range("J22").Value = 2109.2
$ws.Range("K22").Value = 123793.89
$ws.Range("L22").Value = 2109.2
$ws.Range("M22").Value = -123498.89
$ws.Range("N22").Value = -2699.2

$ws.Range("H27").Value = 59749.316
$ws.Range("I27").Value = 123793.89
$ws.Range("J27").Value = 2109.2
$ws.Range("K27").Value = 123793.89
$ws.Range("L27").Value = 2109.2
$ws.Range("M27").Value = -123686.89
$ws.Range("N27").Value = -2323.2

$ws.Range("H46").Value = 14781.477
$ws.Range("I46").Value = 20872.25
$ws.Range("J46").Value = 6660.4443
$ws.Range("K46").Value = 20872.25
$ws.Range("L46").Value = 6660.4443
$ws.Range("M46").Value = -20684.25
$ws.Range("N46").Value = -7036.4443

$ws.Range("H55").Value = 565.7059
$ws.Range("J55").Value = 755.625
$ws.Range("L55").Value = 755.625
$ws.Range("N55").Value = -1101.625

$ws.Range("H93").Value = 3232.6155
$ws.Range("I93").Value = 3466.375
$ws.Range("K93").Value = 3466.375
$ws.Range("M93").Value = -2218.375

$ws.Range("H126").Value = 3999.25
$ws.Range("I126").Value = 3932.3333
$ws.Range("K126").Value = 11796.9999
$ws.Range("M126").Value = -9326.999899999999

$ws.Range("H132").Value = 65051.8
$ws.Range("I132").Value = 79002.25
$ws.Range("K132").Value = 237006.75
$ws.Range("M132").Value = -234476.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 17284.715
$ws.Range("J4").Value = 15332.167
$ws.Range("L4").Value = 15332.167
$ws.Range("N4").Value = -15558.167

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H107").Value = 1527.4546
$ws.Range("I107").Value = 733.8333
$ws.Range("J107").Value = 2479.8
$ws.Range("K107").Value = 2201.4999
$ws.Range("L107").Value = 7439.400000000001
$ws.Range("M107").Value = -281.4998999999998
$ws.Range("N107").Value = -11279.4

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 28729.184
$ws.Range("I126").Value = 31439.059
$ws.Range("K126").Value = 94317.177
$ws.Range("M126").Value = -91847.177
